$wb = $excel.ActiveWorkbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "L6"

# Header row (B1:E1) -- written first so these four strings land at the start of the new shared-string block
$newSheet.Range("B1").Value = "Form"
$newSheet.Range("C1").Value = "Goals scored"
$newSheet.Range("D1").Value = "Goals conceded"
$newSheet.Range("E1").Value = "Total Goals"

# Force column A to text so the "1".."14" rank labels keep their textual shared-string type
$newSheet.Range("A2:A15").NumberFormat = "@"

# Column A: row index labels (written next, to match shared-string allocation order)
$newSheet.Cells.Item(2, 1).Value = "1"
$newSheet.Cells.Item(3, 1).Value = "2"
$newSheet.Cells.Item(4, 1).Value = "3"
$newSheet.Cells.Item(5, 1).Value = "4"
$newSheet.Cells.Item(6, 1).Value = "5"
$newSheet.Cells.Item(7, 1).Value = "6"
$newSheet.Cells.Item(8, 1).Value = "7"
$newSheet.Cells.Item(9, 1).Value = "8"
$newSheet.Cells.Item(10, 1).Value = "9"
$newSheet.Cells.Item(11, 1).Value = "10"
$newSheet.Cells.Item(12, 1).Value = "11"
$newSheet.Cells.Item(13, 1).Value = "12"
$newSheet.Cells.Item(14, 1).Value = "13"
$newSheet.Cells.Item(15, 1).Value = "14"

# Column B: Form data
$newSheet.Range("B2").Value = "AEK,L L W W D L"
$newSheet.Range("B3").Value = "Apollon,W L D D L D"
$newSheet.Range("B4").Value = "Aris,L W L W L D"
$newSheet.Range("B5").Value = "Asteras Tripolis,D L D L L D"
$newSheet.Range("B6").Value = "Atromitos,D L D L D D"
$newSheet.Range("B7").Value = "Giannina,W W L L D L"
$newSheet.Range("B8").Value = "Lamia,D W D W D D"
$newSheet.Range("B9").Value = "Larisa,L W D W D L"
$newSheet.Range("B10").Value = "OFI Crete,L D W D D W"
$newSheet.Range("B11").Value = "Olympiakos,W W W L W W"
$newSheet.Range("B12").Value = "Panathinaikos,D W L L D D"
$newSheet.Range("B13").Value = "Panetolikos,L L D L W W"
$newSheet.Range("B14").Value = "PAOK,W L D W W D"
$newSheet.Range("B15").Value = "Volos NFC,W D D W D D"

# Column C: Goals scored data
$newSheet.Range("C2").Value = "AEK,1 1 3 3 1 0"
$newSheet.Range("C3").Value = "Apollon,1 0 1 0 0 0"
$newSheet.Range("C4").Value = "Aris,0 2 1 2 0 1"
$newSheet.Range("C5").Value = "Asteras Tripolis,2 0 1 1 0 1"
$newSheet.Range("C6").Value = "Atromitos,0 0 1 0 1 0"
$newSheet.Range("C7").Value = "Giannina,1 1 1 1 1 0"
$newSheet.Range("C8").Value = "Lamia,0 3 1 2 0 0"
$newSheet.Range("C9").Value = "Larisa,1 2 1 1 0 0"
$newSheet.Range("C10").Value = "OFI Crete,0 0 2 0 1 1"
$newSheet.Range("C11").Value = "Olympiakos,1 5 3 0 1 2"
$newSheet.Range("C12").Value = "Panathinaikos,2 3 1 1 1 0"
$newSheet.Range("C13").Value = "Panetolikos,0 0 1 1 1 1"
$newSheet.Range("C14").Value = "PAOK,3 0 1 2 1 0"
$newSheet.Range("C15").Value = "Volos NFC,1 0 1 3 1 0"

# Column D: Goals conceded data
$newSheet.Range("D2").Value = "AEK,3 5 1 1 1 2"
$newSheet.Range("D3").Value = "Apollon,0 2 1 0 1 0"
$newSheet.Range("D4").Value = "Aris,1 0 3 1 1 1"
$newSheet.Range("D5").Value = "Asteras Tripolis,2 2 1 3 1 1"
$newSheet.Range("D6").Value = "Atromitos,0 1 1 1 1 0"
$newSheet.Range("D7").Value = "Giannina,0 0 2 2 1 1"
$newSheet.Range("D8").Value = "Lamia,0 0 1 1 0 0"
$newSheet.Range("D9").Value = "Larisa,3 0 1 0 0 1"
$newSheet.Range("D10").Value = "OFI Crete,1 0 1 0 1 0"
$newSheet.Range("D11").Value = "Olympiakos,0 1 1 2 0 0"
$newSheet.Range("D12").Value = "Panathinaikos,2 0 3 2 1 0"
$newSheet.Range("D13").Value = "Panetolikos,1 3 1 3 0 0"
$newSheet.Range("D14").Value = "PAOK,1 3 1 0 0 0"
$newSheet.Range("D15").Value = "Volos NFC,0 0 1 1 1 0"

# Column E: Total Goals data
$newSheet.Range("E2").Value = "AEK,4 6 4 4 2 2"
$newSheet.Range("E3").Value = "Apollon,1 2 2 0 1 0"
$newSheet.Range("E4").Value = "Aris,1 2 4 3 1 2"
$newSheet.Range("E5").Value = "Asteras Tripolis,4 2 2 4 1 2"
$newSheet.Range("E6").Value = "Atromitos,0 1 2 1 2 0"
$newSheet.Range("E7").Value = "Giannina,1 1 3 3 2 1"
$newSheet.Range("E8").Value = "Lamia,0 3 2 3 0 0"
$newSheet.Range("E9").Value = "Larisa,4 2 2 1 0 1"
$newSheet.Range("E10").Value = "OFI Crete,1 0 3 0 2 1"
$newSheet.Range("E11").Value = "Olympiakos,1 6 4 2 1 2"
$newSheet.Range("E12").Value = "Panathinaikos,4 3 4 3 2 0"
$newSheet.Range("E13").Value = "Panetolikos,1 3 2 4 1 1"
$newSheet.Range("E14").Value = "PAOK,4 3 2 2 1 0"
$newSheet.Range("E15").Value = "Volos NFC,1 0 2 4 2 0"

